# Fixing network data cleaning scripts
# - drop the trailing footnote/metadata rows (465-469 and 476-480)
# - rename header row to the cleaned column names
# - title-case the Spanish connector words (de/del/la/las/los/el) in the
#   state/municipality text columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing blocks of footnote rows. Delete the lower block
# first so the row numbers of the upper block don't shift before we get to it.
$ws.Range("A476:D480").EntireRow.Delete()
$ws.Range("A465:D469").EntireRow.Delete()

# Title-case the Spanish connector words wherever they appear as a whole
# word inside column A / B text values (e.g. "Mazapa de Madero" ->
# "Mazapa De Madero").
for ($r = 1; $r -le 463; $r++) {
    foreach ($c in 1, 2) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null -and $v -is [string]) {
            $new = $v -replace '\bde\b', 'De' `
                       -replace '\bdel\b', 'Del' `
                       -replace '\blas\b', 'Las' `
                       -replace '\bla\b', 'La' `
                       -replace '\blos\b', 'Los' `
                       -replace '\bel\b', 'El'
            $cell.Value = $new
        }
    }
}

# Rename the header row to the cleaned/English column names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
